# Generate Report for Handoff
# - Files "5ce00ed3-418d-4574-a190-6e2c07e95f60.md" and
#   "724aa454-5fdd-4e09-9e8d-a0f140b472b3.md" were re-handed-off as a
#   single combined package. They are renamed, their status flips from
#   "Handed back: in sync with en-US" to "Ready for handoff", and the
#   per-locale sheets get fresh handoff file/datetime values while the
#   previous handback file/datetime are cleared (no handback yet).

$wb = $excel.ActiveWorkbook

$oldMd1 = "5ce00ed3-418d-4574-a190-6e2c07e95f60.md"
$oldMd2 = "724aa454-5fdd-4e09-9e8d-a0f140b472b3.md"
$newMd1 = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$newMd2 = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"

$newStatus = "Ready for handoff"

$oldXlfZh1 = "5ce00ed3-418d-4574-a190-6e2c07e95f60.eea2bed3c082eab2d5b0b2f84109bde38c034b4c.zh-cn.xlf"
$oldXlfZh2 = "724aa454-5fdd-4e09-9e8d-a0f140b472b3.621afe6e4bf1e83d4f4b328225316bec895201f1.zh-cn.xlf"
$newXlfZh  = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"

$oldXlfDe1 = "5ce00ed3-418d-4574-a190-6e2c07e95f60.eea2bed3c082eab2d5b0b2f84109bde38c034b4c.de-de.xlf"
$oldXlfDe2 = "724aa454-5fdd-4e09-9e8d-a0f140b472b3.621afe6e4bf1e83d4f4b328225316bec895201f1.de-de.xlf"
$newXlfDe  = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"

$newHandoffDateZh = "2016-03-10 01:15:40"
$newHandoffDateDe = "2016-03-10 01:15:46"
$zeroDate = "0001-01-01 00:00:00"

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/42969b93968724aabb041a067c19873a75c5926a/e2e/"
$handoffZhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1bbd4e3e836736d953f591f3d076e3e724e73d9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$handoffDeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19e275fb263c0d13abcc7f5ec2e4929cab50b206/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# ---------------------------------------------------------------------
# Sheet "Overview": file names and status text are shown per-locale;
# the hyperlink display text also references the old file names.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($repoBase + $newMd1), "", "", $newMd1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($repoBase + $newMd2), "", "", $newMd2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), ($repoBase + ".localization-config"), "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Value = $newXlfZh
$wsZh.Range("D2").Value = $newHandoffDateZh
$wsZh.Range("E2").Clear()
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Value = $zeroDate

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("C3").Value = $newXlfZh
$wsZh.Range("D3").Value = $newHandoffDateZh
$wsZh.Range("E3").Clear()
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Value = $zeroDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($repoBase + $newMd1), "", "", $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), ($handoffZhBase + $newXlfZh), "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($repoBase + $newMd2), "", "", $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), ($handoffZhBase + $newXlfZh), "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ($repoBase + ".localization-config"), "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Value = $newXlfDe
$wsDe.Range("D2").Value = $newHandoffDateDe
$wsDe.Range("E2").Clear()
$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Value = $zeroDate

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("C3").Value = $newXlfDe
$wsDe.Range("D3").Value = $newHandoffDateDe
$wsDe.Range("E3").Clear()
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Value = $zeroDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($repoBase + $newMd1), "", "", $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), ($handoffDeBase + $newXlfDe), "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($repoBase + $newMd2), "", "", $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), ($handoffDeBase + $newXlfDe), "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ($repoBase + ".localization-config"), "", "", ".localization-config")

Write-Output "Done: generated handoff report."
